$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel;
# force a Text number format first so they stay text, like the source data.
$textCells = @("D5", "D6", "D14", "D15", "D19", "D20", "D21", "D23", "D25", "D31", "D33", "D34", "D35", "D38", "D39", "D41", "D43", "D44", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "58.000.53"
$ws.Range("E2").Value = "  +2.88%  "
$ws.Range("D3").Value = "3.065.35"
$ws.Range("E3").Value = "  +3.07%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "527.13"
$ws.Range("E5").Value = "  +6.17%  "
$ws.Range("D6").Value = "142.28"
$ws.Range("E6").Value = "  +5.23%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +5.15%  "
$ws.Range("E9").Value = "  +7.14%  "
$ws.Range("E10").Value = "  +7.79%  "
$ws.Range("E11").Value = "  +5.65%  "
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("D13").Value = "3.588.60"
$ws.Range("E13").Value = "  +3.06%  "
$ws.Range("D14").Value = "27.28"
$ws.Range("E14").Value = "  +8.63%  "
$ws.Range("D15").Value = "0.0000171"
$ws.Range("E15").Value = "  +16.61%  "
$ws.Range("D16").Value = "57.941.28"
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("E17").Value = "  +7.20%  "
$ws.Range("D18").Value = "3.067.07"
$ws.Range("D19").Value = "13.16"
$ws.Range("E19").Value = "  +6.51%  "
$ws.Range("D20").Value = "8.18"
$ws.Range("E20").Value = "  +5.54%  "
$ws.Range("D21").Value = "340.25"
$ws.Range("E21").Value = "  +4.28%  "
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").Value = "5.68"
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D25").Value = "64.99"
$ws.Range("E25").Value = "  +5.67%  "
$ws.Range("E26").Value = "  +6.74%  "
$ws.Range("D27").Value = "0.0₃0980"
$ws.Range("E27").Value = "  +8.83%  "
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("E29").Value = "  +7.39%  "
$ws.Range("E30").Value = "  +8.98%  "
$ws.Range("D31").Value = "1.86"
$ws.Range("E31").Value = "  +6.42%  "
$ws.Range("E32").Value = "  +6.23%  "
$ws.Range("D33").Value = "21.13"
$ws.Range("E33").Value = "  +4.31%  "
$ws.Range("D34").Value = "4.78"
$ws.Range("E34").Value = "  +7.49%  "
$ws.Range("D35").Value = "156.54"
$ws.Range("E35").Value = "  +2.74%  "
$ws.Range("E37").Value = "  +4.24%  "
$ws.Range("D38").Value = "26.54"
$ws.Range("E38").Value = "  +14.14%  "
$ws.Range("D39").Value = "0.0705"
$ws.Range("E39").Value = "  +5.02%  "
$ws.Range("D40").Value = "3.099.85"
$ws.Range("E40").Value = "  +3.14%  "
$ws.Range("D41").Value = "37.87"
$ws.Range("E41").Value = "  +3.43%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.668"
$ws.Range("E43").Value = "  +5.02%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  +6.07%  "
$ws.Range("D46").Value = "2.331.88"
$ws.Range("E46").Value = "  +5.58%  "
$ws.Range("E47").Value = "  +3.60%  "
$ws.Range("E48").Value = "  +3.69%  "
$ws.Range("D49").Value = "0.0246"
$ws.Range("E49").Value = "  +3.45%  "
$ws.Range("E50").Value = "  +5.43%  "
$ws.Range("D51").Value = "20.27"
$ws.Range("E51").Value = "  +6.64%  "

# Restore default number format (value already committed as text, so this
# is purely cosmetic and keeps styling minimal/close to the original).
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
}
